$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: F14 / G14 ---
$ws.Range("F14").Value = "1.5 + 1"
$ws.Range("G14").Value = "24-8-23 & 25-8-23, 1-9-23"

# --- Row 16: F16 / G16 ---
$ws.Range("F16").Value = "1 + 0.5"
$ws.Range("G16").Value = "1/9/2023 & 2/9/23"

# --- Row 15: F15 / G15 ---
$ws.Range("F15").Value = "1 + 0.75"
$ws.Range("G15").Value = "31-8-23 & 2-9-23"

# --- Row 17: F17 / G17 (new numeric entries, G17 is a date) ---
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 44966

# Apply the date number format (maps to built-in numFmtId 14) to the new date cell
# and to G16 which shares the same new style in the source workbook.
$ws.Range("G16").NumberFormat = "mm-dd-yy"
$ws.Range("G17").NumberFormat = "mm-dd-yy"

# --- Column G width change (17.85546875 -> 23 stored width) ---
$ws.Columns("G").ColumnWidth = 22.17

# --- Sheet view / selection update ---
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$ws.Range("G18").Select()

Write-Host "done"
